# Edit script: "filled benchmarks for mandelbrot and more main text"
#
# Logical changes:
#  1. Reorder sheets: move "cluster ut" to sit before "Sharded Matrix"
#  2. Add a new header "2 shards" in K1 of "Sharded Matrix" and widen column J (bestFit)
#  3. Add a brand-new worksheet "Sharded Mandelbrot" at the end, populated with
#     benchmark data (Average/Stdev for 10 shards vs 1 device, 100000 iterations)
#     and ratio formulas, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Reorder: put "cluster ut" before "Sharded Matrix" -------------------
$clusterUt     = $wb.Worksheets.Item("cluster ut")
$shardedMatrix = $wb.Worksheets.Item("Sharded Matrix")
$clusterUt.Move($shardedMatrix)

# --- 3. New worksheet "Sharded Mandelbrot" -----------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$mandelbrot = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$mandelbrot.Name = "Sharded Mandelbrot"
$excel.ActiveWindow.DisplayRuler = $false

# Match the page setup used throughout the rest of the workbook
$ps = $mandelbrot.PageSetup
$ps.LeftMargin   = 54
$ps.RightMargin  = 54
$ps.TopMargin    = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36
$ps.Orientation  = 1

# Sub-headers (row 3) -- entered first, establishing shared-string order
$mandelbrot.Range("B3").Value = "Average"
$mandelbrot.Range("C3").Value = "Stdev"
$mandelbrot.Range("E3").Value = "Average"
$mandelbrot.Range("F3").Value = "Stdev"

# Headers (row 1)
$mandelbrot.Range("A1").Value = "2 Devices 1Gbit"
$mandelbrot.Range("B1").Value = "10 shards"

# --- 2. Sharded Matrix: new column header + column width --------------------
$shardedMatrix.Range("K1").Value = "2 shards"
$shardedMatrix.Columns.Item(10).AutoFit() | Out-Null
$shardedMatrix.Range("K2").Select() | Out-Null

# back to Sharded Mandelbrot headers
$mandelbrot.Range("C1").Value = "100000 iterations"
$mandelbrot.Range("E1").Value = "1 Device"

# Data rows 4-11: A=matrix size, B/C=10 shards avg/stdev, E/F=1 device avg/stdev, H=1 constant
$mandelbrot.Range("A4").Value = 3000
$mandelbrot.Range("B4").Value = 7561.8
$mandelbrot.Range("C4").Value = 31.84
$mandelbrot.Range("E4").Value = 13659
$mandelbrot.Range("F4").Value = 407.35
$mandelbrot.Range("H4").Value = 1

$mandelbrot.Range("A5").Value = 4000
$mandelbrot.Range("B5").Value = 13699.2
$mandelbrot.Range("C5").Value = 914.66
$mandelbrot.Range("E5").Value = 24844.400000000001
$mandelbrot.Range("F5").Value = 925.68
$mandelbrot.Range("H5").Value = 1

$mandelbrot.Range("A6").Value = 5000
$mandelbrot.Range("B6").Value = 20991.8
$mandelbrot.Range("C6").Value = 510.17
$mandelbrot.Range("E6").Value = 39092.199999999997
$mandelbrot.Range("F6").Value = 1001.62
$mandelbrot.Range("H6").Value = 1

$mandelbrot.Range("A7").Value = 6000
$mandelbrot.Range("B7").Value = 30064.6
$mandelbrot.Range("C7").Value = 774.91
$mandelbrot.Range("E7").Value = 56131.199999999997
$mandelbrot.Range("F7").Value = 1591.17
$mandelbrot.Range("H7").Value = 1

$mandelbrot.Range("A8").Value = 7000
$mandelbrot.Range("B8").Value = 40989.199999999997
$mandelbrot.Range("C8").Value = 1026.77
$mandelbrot.Range("E8").Value = 76742.399999999994
$mandelbrot.Range("F8").Value = 1193.58
$mandelbrot.Range("H8").Value = 1

$mandelbrot.Range("A9").Value = 8000
$mandelbrot.Range("B9").Value = 53303
$mandelbrot.Range("C9").Value = 1001.87
$mandelbrot.Range("E9").Value = 100000.6
$mandelbrot.Range("F9").Value = 1337.9
$mandelbrot.Range("H9").Value = 1

$mandelbrot.Range("A10").Value = 9000
$mandelbrot.Range("B10").Value = 67699.199999999997
$mandelbrot.Range("C10").Value = 719.62
$mandelbrot.Range("E10").Value = 126908.2
$mandelbrot.Range("F10").Value = 1108.58
$mandelbrot.Range("H10").Value = 1

$mandelbrot.Range("A11").Value = 10000
$mandelbrot.Range("B11").Value = 83704.2
$mandelbrot.Range("C11").Value = 1197.97
$mandelbrot.Range("E11").Value = 156677
$mandelbrot.Range("F11").Value = 1194.72
$mandelbrot.Range("H11").Value = 1

# Ratio formulas column I: I4 standalone, I5:I11 shared fill-down
$mandelbrot.Range("I4").Formula = "=E4/B4"
$mandelbrot.Range("I5:I11").Formula = "=E5/B5"

$mandelbrot.Range("F11").Select() | Out-Null
